# Artisan "eventbuttons.xlsx" - add Command adjustSV for relative SV adjustments
#
# Net effect described by the target diff:
#  - A new row is inserted into the "Commands" sheet right before the existing
#    "pidSV(<float>)" row (old row 79), containing the new adjustSV command and
#    its description.
#  - The existing "pidSV(<float>)" entry (now shifted one row down) is renamed
#    to "pidSV(<int>)" since its argument type changed from float to int.
#  - Everything below shifts down by one row (handled automatically by the row
#    insert).
#  - Selection / view state is updated to point at the new row.

$wb = $excel.ActiveWorkbook

$wsButtons  = $wb.Worksheets.Item("Buttons")
$wsOptions  = $wb.Worksheets.Item("Options")
$wsCommands = $wb.Worksheets.Item("Commands")

# --- Commands sheet: insert the new adjustSV row -----------------------------
$wsCommands.Rows("79:79").Insert() | Out-Null

$wsCommands.Range("B79").Value = "adjustSV(<int>)"
$wsCommands.Range("C79").Value = "increases or decreases the current target SV value by <int>"

# the former pidSV(<float>) row, now shifted down to row 80, changes its
# argument type from float to int
$wsCommands.Range("B80").Value = "pidSV(<int>)"

# --- View / selection state ---------------------------------------------------
$wsButtons.Activate() | Out-Null
$wsButtons.Range("A1").Select() | Out-Null

$wsOptions.Activate() | Out-Null
$wsOptions.Range("B5").Select() | Out-Null

$wsCommands.Activate() | Out-Null
$wsCommands.Application.ActiveWindow.ScrollRow = 69
$wsCommands.Range("B79:C79").Select() | Out-Null
